$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.692.82"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.606.67"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.49%  "
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.835.82"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "1.597.09"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E14").Value = "  +4.59%  "
$ws.Range("D15").Value = "29.696.72"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.70%  "
$ws.Range("D20").Value = "0.0₃0700"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.108"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("D34").Value = "1.429.30"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.85%  "
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.550"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "57.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("E42").Value = "  +6.66%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.982"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "1.745.16"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("E51").Value = "  +3.39%  "
